$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Formula = "=""37.767.53"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E2").Value = "  -0.68%  "
$cell = $ws.Range("D3")
$cell.Formula = "=""2.033.01"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -0.14%  "
$cell = $ws.Range("D6")
$cell.Formula = "=""0.606"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E6").Value = "  -1.36%  "
$cell = $ws.Range("D7")
$cell.Formula = "=""60.17"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.81%  "
$cell = $ws.Range("D10")
$cell.Formula = "=""0.0824"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("E11").Value = "  +0.36%  "
$cell = $ws.Range("D12")
$cell.Formula = "=""14.66"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E12").Value = "  -0.16%  "
$cell = $ws.Range("D13")
$cell.Formula = "=""2.333.83"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E13").Value = "  -1.07%  "
$cell = $ws.Range("D14")
$cell.Formula = "=""21.06"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E14").Value = "  +0.51%  "
$cell = $ws.Range("D15")
$cell.Formula = "=""0.772"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("E16").Value = "  -2.22%  "
$cell = $ws.Range("D17")
$cell.Formula = "=""2.023.49"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E17").Value = "  -1.32%  "
$cell = $ws.Range("D18")
$cell.Formula = "=""37.717.19"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.73%  "
$cell = $ws.Range("D19")
$cell.Formula = "=""69.61"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E19").Value = "  -0.09%  "
$cell = $ws.Range("D20")
$cell.Formula = "=""5.89"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E20").Value = "  -5.51%  "
$ws.Range("E21").Value = "  -0.36%  "
$cell = $ws.Range("D22")
$cell.Formula = "=""223.85"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  +0.08%  "
$cell = $ws.Range("D24")
$cell.Formula = "=""2.41"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E25").Value = "  +3.25%  "
$cell = $ws.Range("D26")
$cell.Formula = "=""9.36"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E26").Value = "  +2.01%  "
$cell = $ws.Range("D27")
$cell.Formula = "=""167.47"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("E28").Value = "  -1.72%  "
$cell = $ws.Range("D29")
$cell.Formula = "=""18.78"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E29").Value = "  -0.98%  "
$cell = $ws.Range("D30")
$cell.Formula = "=""1.27"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("E32").Value = "  +8.80%  "
$ws.Range("E33").Value = "  -2.62%  "
$cell = $ws.Range("D34")
$cell.Formula = "=""0.0605"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E34").Value = "  +0.59%  "
$cell = $ws.Range("D35")
$cell.Formula = "=""4.51"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("E36").Value = "  +3.45%  "
$cell = $ws.Range("D37")
$cell.Formula = "=""2.32"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("E38").Value = "  +5.07%  "
$ws.Range("E39").Value = "  -0.12%  "
$cell = $ws.Range("D40")
$cell.Formula = "=""18.15"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E40").Value = "  +8.19%  "
$cell = $ws.Range("D41")
$cell.Formula = "=""1.536.04"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("E42").Value = "  +0.27%  "
$cell = $ws.Range("D43")
$cell.Formula = "=""96.19"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("E44").Value = "  -2.54%  "
$cell = $ws.Range("D45")
$cell.Formula = "=""0.0910"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("E46").Value = "  -1.63%  "
$cell = $ws.Range("D47")
$cell.Formula = "=""4.02"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  -0.04%  "
$cell = $ws.Range("D49")
$cell.Formula = "=""2.97"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E49").Value = "  +0.09%  "
$cell = $ws.Range("D50")
$cell.Formula = "=""7.07"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E50").Value = "  +0.65%  "
$cell = $ws.Range("D51")
$cell.Formula = "=""2.223.24"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("E51").Value = "  -1.13%  "
